$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had a two-row header (row 1 + row 2) followed by 11 data
# rows (rows 3-13). The edit consolidates the header into a single row and
# adds five new leading columns (idx, idx2, Name, Date Start, Date End),
# renaming/extending the existing measurement headers. This removes one row
# overall (old row 2), so every data row shifts up by one.

# Remove the old second header row ("Hiver/Eté/Année" sub-header); this also
# shifts data rows 3..13 up to 2..12 and fixes up the sheet dimension
# automatically.
$ws.Rows.Item(2).Delete() | Out-Null

# --- Row 1: consolidated header row -------------------------------------
# Columns A:E are brand-new plain (unstyled) header cells.
$ws.Range("A1:E1").Style = "Normal"
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# Columns F:K keep/receive the measurement-unit header styling (Arial 9,
# general format) that is used throughout the rest of the table.
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.Name = "Arial"

# --- Sheet view / dimension ----------------------------------------------
$ws.Range("A2:K2").Select() | Out-Null
